# Preliminary work so that readInputSheet reads ALL files in the test_files
# directory, not just the sixteen_tests, etc.
#
# On the "optimization_parameters" sheet:
#   - row 8  "Model" -> "production_function" (value stays "MM")
#   - a new row is inserted right after it: "L_curve" / 0
#   - the old "Deletion" row (now pushed down to row 17) is removed
#   - the now-unused trailing duplicate "value" header cells (C1:F1) are cleared
#
# The active sheet/selection also moves from network_optimized_weights!A17
# to optimization_parameters!E9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# "Model" -> "production_function" (row 8, column A)
$ws.Cells.Item(8, 1).Value() = "production_function"

# Insert the new "L_curve" row right after it, shifting everything below down.
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).Value() = "L_curve"
$ws.Cells.Item(9, 2).Value() = 0
$ws.Cells.Item(9, 2).NumberFormat() = "0.00E+00"

# Remove the old "Deletion" row, now sitting at row 17 after the insert above.
$ws.Rows.Item(17).Delete()

# Drop the stray repeated "value" cells in the header row.
$ws.Range("C1:F1").ClearContents()

# Move the active selection: optimization_parameters becomes the active
# sheet/tab, selection moves to E9 (was network_optimized_weights!A17).
$ws.Activate() | Out-Null
$ws.Range("E9").Select() | Out-Null
